# Append 3 new daily rows (2025-12-30, 2025-12-31, 2026-01-01) to the
# "Chart" sheet's GSC export table, matching the export script's usual
# row shape: Date, Impressions-ish counters, and a Reason column that is
# blank for fully-valid rows.
#
# Column A holds dates stored as literal text (e.g. "2025-12-29"), not
# real Excel date serials, so every new date is entered with a leading
# apostrophe to force text entry and stop Excel's automatic date
# recognition from rewriting it as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A88").Value = "'2025-12-30"
$ws.Range("B88").Value = 22
$ws.Range("C88").Value = 1
$ws.Range("D88").Value = 0

$ws.Range("A89").Value = "'2025-12-31"
$ws.Range("B89").Value = 22
$ws.Range("C89").Value = 1
$ws.Range("D89").Value = 0

$ws.Range("A90").Value = "'2026-01-01"
$ws.Range("B90").Value = 22
$ws.Range("C90").Value = 1
# Last row's Reason/validation text column is blank (empty string, not
# an empty/untouched cell) - a lone quote enters an empty text value.
$ws.Range("D90").Value = "'"
